$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ58699681",
    "summ59046682",
    "summ59403291",
    "summ59754824",
    "summ00119994",
    "summ00486599",
    "summ00837041",
    "summ01197525",
    "summ01554889"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
